# Auto-generated script applying the diff to the last paragraph of the document.
$d = $word.ActiveDocument

# Step 1: locate the old full sentence (unique in the document) and shorten it
# to just the first replacement run's text, keeping it as a single run with no
# special run properties (matches the first <w:r><w:t>s we can see</w:t></w:r>).
$oldSentence = 's we can see, for the topic 1, the key words in the news headline is covid.'
$findRange = $d.Content
$found = $findRange.Find.Execute($oldSentence, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "could not find the original sentence to edit" }
$findRange.Text = 's we can see'

# Step 2: get a fresh (non-collapsed) anchor range on that now-unique snippet; each
# subsequent InsertXML call appends a new run immediately after the current content
# of this paragraph, in order, without disturbing anything already written.
$anchor = $d.Content
$found2 = $anchor.Find.Execute('s we can see', $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "could not find the anchor snippet" }

function Add-Run($rngObj, $text, $hint, $preserve, $lastBreak) {
    $rPrInner = ""
    if ($hint) { $rPrInner += '<w:rFonts w:hint="eastAsia"/>' }
    $rPr = ""
    if ($rPrInner -ne "") { $rPr = "<w:rPr>" + $rPrInner + "</w:rPr>" }
    $breakTag = ""
    if ($lastBreak) { $breakTag = "<w:lastRenderedPageBreak/>" }
    $spaceAttr = ""
    if ($preserve) { $spaceAttr = ' xml:space="preserve"' }
    $escapedText = $text -replace "&", "&amp;" -replace "<", "&lt;" -replace ">", "&gt;"
    $xml = '<?xml version="1.0" encoding="utf-8"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r>' + $rPr + $breakTag + '<w:t' + $spaceAttr + '>' + $escapedText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rngObj.InsertXML($xml)
}

# Step 3: append every remaining run, in order, onto the anchor range.
Add-Run $anchor ' from the result, the words “Australian, say, policy' $false $true $false
Add-Run $anchor ', men, women' $false $false $false
Add-Run $anchor '” ' $true $true $false
Add-Run $anchor 'appear a lot of times, which is reasonable, since our dataset comes from ' $false $true $false
Add-Run $anchor 'Australian Broadcasting Corporation' $false $false $false
Add-Run $anchor ' and ' $false $true $false
Add-Run $anchor 'other' $false $false $false
Add-Run $anchor ' words are the common words used in the news. Ignore ' $false $true $false
Add-Run $anchor 'these words' $false $false $false
Add-Run $anchor ',' $false $false $false
Add-Run $anchor ' the key words in the news headline ' $false $true $false
Add-Run $anchor 'in 2017 is “' $false $false $false
Add-Run $anchor 't' $false $false $false
Add-Run $anchor 'rump' $true $false $false
Add-Run $anchor '”' $false $false $false
Add-Run $anchor ' “government”' $false $true $false
Add-Run $anchor ' ' $true $true $false
Add-Run $anchor '“' $true $false $false
Add-Run $anchor 'fire' $false $false $false
Add-Run $anchor '”;' $true $false $false
Add-Run $anchor ' ' $false $true $false
Add-Run $anchor ' ' $true $true $false
Add-Run $anchor 'in' $true $false $false
Add-Run $anchor ' 2018 ' $false $true $false
Add-Run $anchor 'are' $true $false $false
Add-Run $anchor ' “crash” “fire” “death”, in 2019 are “fire” “murder” “bushfire”;  in 2020 are “covid” “election” ' $false $true $false
Add-Run $anchor '“restriction”; in 2021 are “covid' $false $false $true
Add-Run $anchor '” “' $false $false $false
Add-Run $anchor 'vaccine' $false $false $false
Add-Run $anchor '” ' $false $true $false
Add-Run $anchor '“' $true $false $false
Add-Run $anchor 'lockdown”' $false $false $false
Add-Run $anchor '.' $true $false $false
Add-Run $anchor ' From these key words we can recall some important event during that year, for instance, in ' $false $true $false
Add-Run $anchor '2017' $false $false $false
Add-Run $anchor ', ' $false $true $false
Add-Run $anchor 'the second year of Trump''s presidency of the United States, he adjusted and promulgated many bills and bans, which attracted the attention of the whole world' $false $false $false
Add-Run $anchor ' and in 2019, t' $false $true $false
Add-Run $anchor 'he coronavirus was just discovered in China and didn''t get the world''s attention' $false $false $false
Add-Run $anchor ', while in the next 2 years, p' $false $false $false
Add-Run $anchor 'eople are forced to lock down because of the coronavirus epidemic' $false $false $false
Add-Run $anchor '. ' $false $true $false

Write-Output "edit complete"
